$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 12.5584
$ws.Range("B4").Value = 5.727699999999994
$ws.Range("C4").Value = -14.18229999999999
$ws.Range("E4").Value = 12.99869999999999
$ws.Range("C5").Value = -14.61310000000001
$ws.Range("B6").Value = 9.403499999999996
$ws.Range("B7").Value = 7.037599999999998
$ws.Range("C8").Value = -12.2213
$ws.Range("E9").Value = 14.49760000000001
$ws.Range("E11").Value = 13.20219999999999
$ws.Range("E14").Value = 13.13570000000001
$ws.Range("B16").Value = 9.041600000000006
$ws.Range("C16").Value = -12.0968
$ws.Range("E18").Value = 12.8799
$ws.Range("B20").Value = 6.102399999999998
$ws.Range("C22").Value = -11.0064
$ws.Range("E25").Value = 12.90679999999999
